$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.923.13"
$ws.Range("E2").Value = "  +4.61%  "

$ws.Range("D3").Value = "'2.231.77"
$ws.Range("E3").Value = "  +4.63%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'251.18"
$ws.Range("E5").Value = "  +7.07%  "

$ws.Range("D6").Value = "'0.613"
$ws.Range("E6").Value = "  +2.87%  "

$ws.Range("D7").Value = "'74.92"
$ws.Range("E7").Value = "  +9.07%  "

$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").Value = "'0.594"
$ws.Range("E9").Value = "  +5.80%  "

$ws.Range("D10").Value = "'41.19"
$ws.Range("E10").Value = "  +8.23%  "

$ws.Range("D11").Value = "'0.0921"
$ws.Range("E11").Value = "  +4.47%  "

$ws.Range("D12").Value = "'6.86"
$ws.Range("E12").Value = "  +5.18%  "

$ws.Range("D13").Value = "'0.102"
$ws.Range("E13").Value = "  +2.80%  "

$ws.Range("D14").Value = "'2.566.98"
$ws.Range("E14").Value = "  +4.56%  "

$ws.Range("D15").Value = "'14.46"
$ws.Range("E15").Value = "  +2.71%  "

$ws.Range("D16").Value = "'2.235.64"
$ws.Range("E16").Value = "  +5.59%  "

$ws.Range("D17").Value = "'0.785"
$ws.Range("E17").Value = "  +2.06%  "

$ws.Range("D18").Value = "'42.823.31"
$ws.Range("E18").Value = "  +4.87%  "

$ws.Range("E19").Value = "  +5.15%  "

$ws.Range("D20").Value = "'71.12"
$ws.Range("E20").Value = "  +3.90%  "

$ws.Range("D21").Value = "'5.93"
$ws.Range("E21").Value = "  +4.80%  "

$ws.Range("D22").Value = "'227.75"
$ws.Range("E22").Value = "  +1.92%  "

$ws.Range("D23").Value = "'2.17"
$ws.Range("E23").Value = "  +13.41%  "

$ws.Range("D24").Value = "'9.40"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").Value = "'10.65"
$ws.Range("E26").Value = "  +2.35%  "

$ws.Range("D27").Value = "'3.42"
$ws.Range("E27").Value = "  +5.09%  "

$ws.Range("D28").Value = "'39.01"
$ws.Range("E28").Value = "  +27.28%  "

$ws.Range("D29").Value = "'2.22"
$ws.Range("E29").Value = "  +4.01%  "

$ws.Range("D30").Value = "'2.22"
$ws.Range("E30").Value = "  +3.15%  "

$ws.Range("D31").Value = "'171.50"
$ws.Range("E31").Value = "  +0.94%  "

$ws.Range("D32").Value = "'20.13"
$ws.Range("E32").Value = "  +3.57%  "

$ws.Range("D33").Value = "'0.0792"
$ws.Range("E33").Value = "  +6.67%  "

$ws.Range("D34").Value = "'5.20"
$ws.Range("E34").Value = "  +3.94%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.112"
$ws.Range("E35").Value = "  +10.27%  "

$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").Value = "'0.121"
$ws.Range("E36").Value = "  +1.97%  "

$ws.Range("D37").Value = "'4.42"
$ws.Range("E37").Value = "  +8.03%  "

$ws.Range("D38").Value = "'0.0325"
$ws.Range("E38").Value = "  +14.37%  "

$ws.Range("D39").Value = "'12.46"
$ws.Range("E39").Value = "  +8.38%  "

$ws.Range("D40").Value = "'2.09"
$ws.Range("E40").Value = "  +3.77%  "

$ws.Range("D41").Value = "'0.204"
$ws.Range("E41").Value = "  +10.73%  "

$ws.Range("D42").Value = "'5.36"
$ws.Range("E42").Value = "  +3.59%  "

$ws.Range("D43").Value = "'59.38"
$ws.Range("E43").Value = "  +4.27%  "

$ws.Range("D44").Value = "'8.64"
$ws.Range("E44").Value = "  +6.47%  "

$ws.Range("B45").Value = "WOONetwork"
$ws.Range("C45").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D45").Value = "'0.482"
$ws.Range("E45").Value = "  +32.19%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'103.06"
$ws.Range("E46").Value = "  +7.38%  "

$ws.Range("D47").Value = "'0.0986"
$ws.Range("E47").Value = "  +3.81%  "

$ws.Range("D48").Value = "'2.42"
$ws.Range("E48").Value = "  +15.13%  "

$ws.Range("E49").Value = "  +4.58%  "

$ws.Range("E50").Value = "  +4.22%  "

$ws.Range("D51").Value = "'2.67"
$ws.Range("E51").Value = "  +2.43%  "
